$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching style of existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-28
$iValues = @(7,6,4,5,5,7,1,1,4,9,6,9,10,1,3,7,9,6,1,5,1,9,1,1,1,1,8)
$jValues = @(7,7,4,6,8,8,5,4,7,9,7,9,10,3,8,8,9,8,5,7,5,9,5,5,5,2,8)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $rowNum = $r + 2
    $ws.Cells.Item($rowNum, 9).Value = $iValues[$r]
    $ws.Cells.Item($rowNum, 10).Value = $jValues[$r]
}
